$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name order (Angola moves up, before Mauritania) ---
# Row 118 becomes Angola, Row 119 becomes Mauritania, Row 120 becomes Lituania
$ws.Range("A118").Value = "Angola"
$ws.Range("A119").Value = "Mauritania"
$ws.Range("A120").Value = "Lituania"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 00:17"

# --- Update numeric statistics cells ---
# Row 4
$ws.Range("B4").Value = 8385541
$ws.Range("C4").Value = 42876
$ws.Range("D4").Value = 5452751
$ws.Range("E4").Value = 2708064
$ws.Range("G4").Value = 444
$ws.Range("H4").Value = 224726

# Row 5
$ws.Range("B5").Value = 7547762
$ws.Range("C5").Value = 55035
$ws.Range("E5").Value = 774196

# Row 6
$ws.Range("B6").Value = 5235344
$ws.Range("C6").Value = 10982
$ws.Range("E6").Value = 446124
$ws.Range("G6").Value = 215
$ws.Range("H6").Value = 153905

# Row 10
$ws.Range("B10").Value = 959572
$ws.Range("C10").Value = 7201
$ws.Range("D10").Value = 858294
$ws.Range("E10").Value = 72308
$ws.Range("G10").Value = 167
$ws.Range("H10").Value = 28970

# Row 21
$ws.Range("B21").Value = 366981
$ws.Range("C21").Value = 5248
$ws.Range("E21").Value = 67115

# Row 31
$ws.Range("B31").Value = 198148
$ws.Range("C31").Value = 1827
$ws.Range("D31").Value = 167112
$ws.Range("E31").Value = 21276

# Row 71
$ws.Range("B71").Value = 47310
$ws.Range("C71").Value = 111
$ws.Range("D71").Value = 46618
$ws.Range("E71").Value = 382

# Row 85
$ws.Range("B85").Value = 29503
$ws.Range("C85").Value = 395
$ws.Range("D85").Value = 16943
$ws.Range("E85").Value = 11574
$ws.Range("G85").Value = 18
$ws.Range("H85").Value = 986

# Row 113
$ws.Range("B113").Value = 8964
$ws.Range("C113").Value = 8
$ws.Range("E113").Value = 1520

# Row 117
$ws.Range("D117").Value = 6526
$ws.Range("E117").Value = 1140

# Row 118
$ws.Range("B118").Value = 7622
$ws.Range("C118").Value = 160
$ws.Range("D118").Value = 3030
$ws.Range("E118").Value = 4345
$ws.Range("G118").Value = 6
$ws.Range("H118").Value = 247

# Row 119
$ws.Range("B119").Value = 7608
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 7347
$ws.Range("E119").Value = 98
$ws.Range("H119").Value = 163

# Row 120
$ws.Range("B120").Value = 7521
$ws.Range("C120").Value = 252
$ws.Range("D120").Value = 3097
$ws.Range("E120").Value = 4311
$ws.Range("H120").Value = 113

# Row 150
$ws.Range("B150").Value = 3388
$ws.Range("C150").Value = 9
$ws.Range("D150").Value = 2586
$ws.Range("E150").Value = 670

# Row 158
$ws.Range("B158").Value = 2381
$ws.Range("C158").Value = 38
$ws.Range("D158").Value = 1774
$ws.Range("E158").Value = 542

# Row 167
$ws.Range("B167").Value = 1210
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 15

# Row 172
$ws.Range("B172").Value = 744
$ws.Range("C172").Value = 29
$ws.Range("D172").Value = 433
$ws.Range("E172").Value = 310

# Row 192
$ws.Range("B192").Value = 221
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 200
$ws.Range("E192").Value = 14
